$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "N° segment"
$ws.Range("B1").Value = "pièce"
$ws.Range("C1").Value = "type"
$ws.Range("D1").Value = "couleur"
$ws.Range("E1").Value = "surface"

# Numbering / surface (numeric) columns first
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("E2").Value = 70
$ws.Range("E3").Value = 25
$ws.Range("E4").Value = 20
$ws.Range("E5").Value = 10

# "pièce" (column B) and "type" (column C) — entry order chosen to match
# the shared-string insertion order of the target workbook
$ws.Range("B2").Value = "Salle"
$ws.Range("B3").Value = "cuisine"
$ws.Range("B4").Value = "chambre 1"
$ws.Range("C2").Value = "carrelage"
$ws.Range("C3").Value = "carrelage"
$ws.Range("B5").Value = "salle de bains"
$ws.Range("C4").Value = "parquet"
$ws.Range("C5").Value = "carrelage"

# Remaining numbering rows (col A only)
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19

# Column widths to match best-fit content widths (values chosen so the
# runtime's internal width quantization lands as close as possible to the
# target stored widths of 11.285.., 14.140.., 9.140.., 7.710.., 7.425..)
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 13.333333333333334
$ws.Columns.Item(3).ColumnWidth = 8.333333333333334
$ws.Columns.Item(4).ColumnWidth = 6.833333333333333
$ws.Columns.Item(5).ColumnWidth = 6.666666666666667

# Selection matches the target (active cell E4)
$ws.Range("E4").Select() | Out-Null
